$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.676.35'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '3.525.61'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''621.58'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('D6').Value = '''171.75'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('D7').Value = '3.517.07'
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').Value = '''7.19'
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '''46.19'
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').Value = '4.083.89'
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '''607.31'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.515.87'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '70.786.45'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').Value = '''17.68'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').Value = '''15.56'
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('D25').Value = '''97.45'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').Value = '''33.62'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '''6.79'
$ws.Range('E34').Value = '  -4.70%  '
$ws.Range('D35').Value = '''615.34'
$ws.Range('E35').Value = '  -5.80%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.0492'
$ws.Range('E36').Value = '  +3.64%  '
$ws.Range('B37').Value = 'Cosmos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D37').Value = '''10.85'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('D39').Value = '''56.70'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = '''3.39'
$ws.Range('E41').Value = '  -6.77%  '
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').Value = '3.339.83'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').Value = '0.0₃0723'
$ws.Range('E44').Value = '  +2.14%  '
$ws.Range('E45').Value = '  -3.30%  '
$ws.Range('D46').Value = '''2.91'
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').Value = '''31.80'
$ws.Range('E47').Value = '  -2.94%  '
$ws.Range('E48').Value = '  -5.42%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '''133.81'
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('B51').Value = 'USDe'
$ws.Range('C51').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D51').Value = '''1.00'
$ws.Range('E51').Value = '  -0.02%  '
